# Clean the data and rewrite the sheet with new headers/columns/values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - A1 is a brand-new header cell; clone B1's header style (s="1")
# onto it via copy/paste-special (format only) before setting its value.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1").Value = "City"
$ws.Range("B1").Value = "Last"
$ws.Range("C1").Value = "State"
$ws.Range("D1").Value = "Salary"
$ws.Range("E1").Value = "Tax %"
$ws.Range("F1").Value = "Taxes Owed"
$ws.Range("G1").Value = "Sal After Tax"
$ws.Range("H1").Value = "Conditional"

# Data rows: City, Last, State, Salary, Tax%, Taxes Owed, Sal After Tax, Conditional
$data = @(
    @("Riverside",   "Doe",       " NJ", 45000,  0.2,  9000,  36000, $true),
    @("Phila",       "McGinnis",  " PA", 18000,  0.15, 2700,  15300, $true),
    @("Riverside",   "Repici",    " NJ", 120000, 0.25, 30000, 90000, $false),
    @("SomeTown",    "Tyler",     "SD",  90000,  0.25, 22500, 67500, $false),
    @("SomeTown",    "Blankman",  " SD", 30000,  0.15, 4500,  25500, $true),
    @("Desert City", "Jet",       "CO",  68000,  0.2,  13600, 54400, $false)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}
